# Generate Report for Handback
#
# This script reproduces, via Excel COM-interop calls, the "handback" report
# generation that:
#   1. Marks every row's Status as handed-back (was "Ready for handoff").
#   2. Fills in the "Latest Target File" (hyperlink) and "Latest Handback
#      File" columns on the zh-cn / de-de sheets with the generated xliff.
#   3. Stamps the "Latest Handback DateTime" for each locale.
#   4. Grows the now-wider columns to fit their new contents.

$wb  = $excel.ActiveWorkbook
$ovw = $wb.Worksheets.Item("Overview")
$zh  = $wb.Worksheets.Item("zh-cn")
$de  = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus
$ovw.Range("E3").Value = $newStatus
$ovw.Range("F3").Value = $newStatus

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. Latest Target File (I, hyperlink to the source .md on GitHub) and
#    Latest Handback File (J, generated xliff name) for zh-cn and de-de.
# ---------------------------------------------------------------------
$mdDisplay = "51b2560d-38fa-4c78-9d91-d689fe473d7d.md"
$mdAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c25a613832dda3dec2fc0492b8e64d192b751808/e2e/51b2560d-38fa-4c78-9d91-d689fe473d7d.md"

$zh.Hyperlinks.Add($zh.Range("I2"), $mdAddress, "", "", $mdDisplay)
$zh.Hyperlinks.Add($zh.Range("I3"), $mdAddress, "", "", $mdDisplay)

$de.Hyperlinks.Add($de.Range("I2"), $mdAddress, "", "", $mdDisplay)
$de.Hyperlinks.Add($de.Range("I3"), $mdAddress, "", "", $mdDisplay)

$zh.Range("J2").Value = "51b2560d-38fa-4c78-9d91-d689fe473d7d.f8bb77d93af7183a36c267851502e5eb24c2419e.zh-cn.xlf"
$zh.Range("J3").Value = "51b2560d-38fa-4c78-9d91-d689fe473d7d.f8bb77d93af7183a36c267851502e5eb24c2419e.zh-cn.xlf"

$de.Range("J2").Value = "51b2560d-38fa-4c78-9d91-d689fe473d7d.f8bb77d93af7183a36c267851502e5eb24c2419e.de-de.xlf"
$de.Range("J3").Value = "51b2560d-38fa-4c78-9d91-d689fe473d7d.f8bb77d93af7183a36c267851502e5eb24c2419e.de-de.xlf"

# ---------------------------------------------------------------------
# 3. Latest Handback DateTime (K) per locale.
# ---------------------------------------------------------------------
$zh.Range("K2").Value = "2016-08-23 11:04:29"
$zh.Range("K3").Value = "2016-08-23 11:04:29"

$de.Range("K2").Value = "2016-08-23 11:04:37"
$de.Range("K3").Value = "2016-08-23 11:04:37"

# ---------------------------------------------------------------------
# 4. Column widths grow now that the Status / Target / Handback columns
#    hold longer text.
# ---------------------------------------------------------------------
$ovw.Range("E:F").ColumnWidth = 29.166666666666668

$zh.Range("C:C").ColumnWidth = 29.166666666666668
$zh.Range("I:J").ColumnWidth = 39.166666666666664

$de.Range("C:C").ColumnWidth = 29.166666666666668
$de.Range("I:J").ColumnWidth = 39.166666666666664
